$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-08-16 Saturday"; New = "2025-08-17 Sunday" },
    @{ Old = "307÷3="; New = "598÷4=" },
    @{ Old = "190÷7="; New = "561÷6=" },
    @{ Old = "505÷4="; New = "505÷8=" },
    @{ Old = "392÷6="; New = "587÷9=" },
    @{ Old = "884÷5="; New = "134÷3=" },
    @{ Old = "583÷5="; New = "267÷4=" },
    @{ Old = "290÷4="; New = "291÷6=" },
    @{ Old = "762÷8="; New = "932÷6=" },
    @{ Old = "976÷5="; New = "239÷9=" },
    @{ Old = "401÷3="; New = "900÷2=" },
    @{ Old = "572÷3="; New = "917÷4=" },
    @{ Old = "389÷3="; New = "794÷8=" },
    @{ Old = "379÷8="; New = "328÷9=" },
    @{ Old = "792÷2="; New = "604÷4=" },
    @{ Old = "562÷4="; New = "843÷7=" },
    @{ Old = "347÷9="; New = "589÷6=" },
    @{ Old = "420÷2="; New = "607÷8=" },
    @{ Old = "114÷3="; New = "919÷3=" },
    @{ Old = "192÷7="; New = "729÷2=" },
    @{ Old = "981÷6="; New = "727÷9=" },
    @{ Old = "891÷4="; New = "488÷8=" },
    @{ Old = "820÷9="; New = "918÷2=" },
    @{ Old = "954÷6="; New = "926÷3=" },
    @{ Old = "239÷5="; New = "577÷3=" },
    @{ Old = "874÷4="; New = "274÷7=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
